# Add a new Leave Card (SL) entry batch to the "2018 LEAVE CREDITS" sheet
# and refresh the active-cell/selection view state to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2018 LEAVE CREDITS")

# --- Row 27: fill in the EARNED value for the existing SL(1-0-0) entry dated 45108 ---
$ws.Range("C27").Value = 1.25

# --- Row 29: new period date + EARNED value ---
$ws.Range("A29").Value = 45139
$ws.Range("C29").Value = 1.25

# --- Row 30: new period date, particulars, EARNED value, absence flag, remarks date ---
$ws.Range("A30").Value = 45170
$ws.Range("B30").Value = "SL(1-0-0)"
$ws.Range("C30").Value = 1.25
$ws.Range("H30").Value = 1
# K30 picks up the same short-date format already used by K27/K28 (the
# "REMARKS" date column) -- copy formats from K27 so the existing style is
# reused instead of a duplicate being created.
$ws.Range("K27").Copy()
$ws.Range("K30").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("K30").Value = 45191

# --- Row 31: new period date + EARNED value ---
$ws.Range("A31").Value = 45230
$ws.Range("C31").Value = 1.25

# --- Refresh the saved view state (scroll/split position + active selection) ---
$ws.Activate()
$win = $ws.Application.ActiveWindow
$win.SplitRow = 12
$ws.Range("A32").Select()
